# Mumbai_India.xlsx — "get batting team and bowling team implemented"
# Fill in the Runs (column B) for the batters now that the batting/bowling
# team data is wired up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1   # KL Rahul
$ws.Range("B3").Value = 4   # Rohit Sharma
$ws.Range("B4").Value = 7   # Virat Kohli(C)
$ws.Range("B5").Value = 8   # Suryakumar Yadav
$ws.Range("B6").Value = 8   # Rishabh Pant
$ws.Range("B11").Value = 7  # Jasprit Bumrah
$ws.Range("B12").Value = 7  # Kuldeep

# Match the author's final on-screen selection/scroll position.
$ws.Range("A2").Select()
$ws.Range("G5").Select()
